# Update "nombre_aides" (C) and "montant_total" (E) figures for the
# 2022-06-14 data refresh of the Fonds de solidarite dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 2;   C = 766330;   E = 1429232685 },
    @{ Row = 13;  C = 187858;   E = 1168030789 },
    @{ Row = 16;  C = 10172;    E = 28006957 },
    @{ Row = 27;  C = 90066;    E = 442910133 },
    @{ Row = 69;  C = 17892;    E = 103956993 },
    @{ Row = 78;  C = 178445;   E = 892645000 },
    @{ Row = 91;  C = 18876;    E = 75305211 },
    @{ Row = 121; C = 1306355;  E = 2275357834 },
    @{ Row = 129; C = 633721;   E = 3433665537 },
    @{ Row = 132; C = 585963;   E = 3471020874 },
    @{ Row = 154; C = 18465;    E = 74174082 },
    @{ Row = 237; C = 283323;   E = 1438433147 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
